$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").ClearFormats()
$ws.Range("E4").Value = 661.68399999999997
$ws.Range("E5").ClearFormats()
$ws.Range("E5").Value = 1800.058
$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = 1333.99
$ws.Range("E7").ClearFormats()
$ws.Range("E7").Value = 1800.088
$ws.Range("E8").ClearFormats()
$ws.Range("E8").Value = 968.20799999999997
$ws.Range("E9").ClearFormats()
$ws.Range("E9").Value = 1098.077
$ws.Range("E10").ClearFormats()
$ws.Range("E10").Value = 1365.8420000000001
$ws.Range("E11").ClearFormats()
$ws.Range("E11").Value = 920.60599999999999
$ws.Range("E12").ClearFormats()
$ws.Range("E12").Value = 1800.0650000000001
$ws.Range("E13").ClearFormats()
$ws.Range("E13").Value = 1800.0619999999999
$ws.Range("E14").ClearFormats()
$ws.Range("E14").Value = 1800.1010000000001
$ws.Range("E15").ClearFormats()
$ws.Range("E15").Value = 1800.183
$ws.Range("E16").ClearFormats()
$ws.Range("E16").Value = 1146.4459999999999
$ws.Range("E17").ClearFormats()
$ws.Range("E17").Value = 1799.0619999999999
$ws.Range("E18").ClearFormats()
$ws.Range("E18").Value = 263.14699999999999
$ws.Range("E19").ClearFormats()
$ws.Range("E19").Value = 1398.829
$ws.Range("E20").ClearFormats()
$ws.Range("E20").Value = 1800.1579999999999
$ws.Range("E21").ClearFormats()
$ws.Range("E21").Value = 750.06999999999994
$ws.Range("E22").ClearFormats()
$ws.Range("E22").Value = 1559.2940000000001
$ws.Range("E23").ClearFormats()
$ws.Range("E23").Value = 856.78600000000006
$ws.Range("E24").ClearFormats()
$ws.Range("E24").Value = 1805.921
$ws.Range("E25").ClearFormats()
$ws.Range("E25").Value = 1817.8969999999999
$ws.Range("E26").ClearFormats()
$ws.Range("E26").Value = 2391.902
$ws.Range("E27").ClearFormats()
$ws.Range("E27").Value = 1819.7860000000001
$ws.Range("E28").ClearFormats()
$ws.Range("E28").Value = 1801.4590000000001
$ws.Range("E33").ClearFormats()
$ws.Range("E33").Value = 720.71600000000001
$ws.Range("E34").ClearFormats()
$ws.Range("E34").Value = 1800.2629999999999
$ws.Range("E35").ClearFormats()
$ws.Range("E35").Value = 1789.6969999999999
$ws.Range("E36").ClearFormats()
$ws.Range("E36").Value = 1800.114
$ws.Range("E37").ClearFormats()
$ws.Range("E37").Value = 900.505
$ws.Range("E38").ClearFormats()
$ws.Range("E38").Value = 1002.476
$ws.Range("E39").ClearFormats()
$ws.Range("E39").Value = 1544.1780000000001
$ws.Range("E40").ClearFormats()
$ws.Range("E40").Value = 1068.412
$ws.Range("E41").ClearFormats()
$ws.Range("E41").Value = 1800.097
$ws.Range("E42").ClearFormats()
$ws.Range("E42").Value = 1800.1489999999999
$ws.Range("E43").ClearFormats()
$ws.Range("E43").Value = 1800.116
$ws.Range("E44").ClearFormats()
$ws.Range("E44").Value = 1800.1120000000001
$ws.Range("E45").ClearFormats()
$ws.Range("E45").Value = 1020.51
$ws.Range("E46").ClearFormats()
$ws.Range("E46").Value = 1800.1869999999999
$ws.Range("E47").ClearFormats()
$ws.Range("E47").Value = 297.95
$ws.Range("E48").ClearFormats()
$ws.Range("E48").Value = 1694.2429999999999
$ws.Range("E49").ClearFormats()
$ws.Range("E49").Value = 1747.6769999999999
$ws.Range("E50").ClearFormats()
$ws.Range("E50").Value = 807.17599999999993
$ws.Range("E51").ClearFormats()
$ws.Range("E51").Value = 1048.952
$ws.Range("E52").ClearFormats()
$ws.Range("E52").Value = 573.03300000000002
$ws.Range("E53").ClearFormats()
$ws.Range("E53").Value = 1803.819
$ws.Range("E54").ClearFormats()
$ws.Range("E54").Value = 1878.6859999999999
$ws.Range("E55").ClearFormats()
$ws.Range("E55").Value = 2006.992
$ws.Range("E56").ClearFormats()
$ws.Range("E56").Value = 1815.8620000000001
$ws.Range("E57").ClearFormats()
$ws.Range("E57").Value = 1800.4680000000001
$ws.Range("E62").ClearFormats()
$ws.Range("E62").Value = 644.779
$ws.Range("E63").ClearFormats()
$ws.Range("E63").Value = 1800.088
$ws.Range("E64").ClearFormats()
$ws.Range("E64").Value = 1800.126
$ws.Range("E65").ClearFormats()
$ws.Range("E65").Value = 1797.8440000000001
$ws.Range("E66").ClearFormats()
$ws.Range("E66").Value = 660.82
$ws.Range("E67").ClearFormats()
$ws.Range("E67").Value = 827.33299999999997
$ws.Range("E68").ClearFormats()
$ws.Range("E68").Value = 943.18799999999999
$ws.Range("E69").ClearFormats()
$ws.Range("E69").Value = 745.553
$ws.Range("E70").ClearFormats()
$ws.Range("E70").Value = 1800.085
$ws.Range("E71").ClearFormats()
$ws.Range("E71").Value = 1759.049
$ws.Range("E72").ClearFormats()
$ws.Range("E72").Value = 1800.07
$ws.Range("E73").ClearFormats()
$ws.Range("E73").Value = 1800.0709999999999
$ws.Range("E74").ClearFormats()
$ws.Range("E74").Value = 951.29700000000003
$ws.Range("E75").ClearFormats()
$ws.Range("E75").Value = 1800.1569999999999
$ws.Range("E76").ClearFormats()
$ws.Range("E76").Value = 254.387
$ws.Range("E77").ClearFormats()
$ws.Range("E77").Value = 1688.0139999999999
$ws.Range("E78").ClearFormats()
$ws.Range("E78").Value = 1800.049
$ws.Range("E79").ClearFormats()
$ws.Range("E79").Value = 1129.0260000000001
$ws.Range("E80").ClearFormats()
$ws.Range("E80").Value = 1438.328
$ws.Range("E81").ClearFormats()
$ws.Range("E81").Value = 784.53700000000003
$ws.Range("E82").ClearFormats()
$ws.Range("E82").Value = 1802.3130000000001
$ws.Range("E83").ClearFormats()
$ws.Range("E83").Value = 1804.8
$ws.Range("E84").ClearFormats()
$ws.Range("E84").Value = 3819.1570000000002
$ws.Range("E85").ClearFormats()
$ws.Range("E85").Value = 3002.4110000000001
$ws.Range("E86").ClearFormats()
$ws.Range("E86").Value = 1800.9960000000001

[void]$ws.Range("L4:L28").Select()
